$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") "63.606.68"
$ws.Range("E2").Value = "  -2.54%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.319.60"
$ws.Range("E3").Value = "  -4.07%  "

# Row 4
$ws.Range("E4").Value = "  +0.50%  "

# Row 5
Set-TextValue $ws.Range("D5") "546.84"
$ws.Range("E5").Value = "  -1.26%  "

# Row 6
Set-TextValue $ws.Range("D6") "171.64"
$ws.Range("E6").Value = "  -4.14%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.611"
$ws.Range("E7").Value = "  -4.01%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D8") "1.00"
$ws.Range("E8").Value = "  +0.11%  "

# Row 9
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextValue $ws.Range("D9") "3.315.94"
$ws.Range("E9").Value = "  -3.57%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.609"
$ws.Range("E10").Value = "  -4.26%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.151"
$ws.Range("E11").Value = "  -1.14%  "

# Row 12
Set-TextValue $ws.Range("D12") "53.30"
$ws.Range("E12").Value = "  -1.53%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000264"
$ws.Range("E13").Value = "  -2.74%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D14") "4.000.16"
$ws.Range("E14").Value = "  +0.61%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D15") "8.84"
$ws.Range("E15").Value = "  -4.52%  "

# Row 16
Set-TextValue $ws.Range("D16") "18.16"
$ws.Range("E16").Value = "  -2.45%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "3.322.20"
$ws.Range("E17").Value = "  -3.28%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D18") "0.116"
$ws.Range("E18").Value = "  -3.62%  "

# Row 19
Set-TextValue $ws.Range("D19") "11.63"
$ws.Range("E19").Value = "  -2.75%  "

# Row 20
Set-TextValue $ws.Range("D20") "63.573.68"
$ws.Range("E20").Value = "  -2.25%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.972"
$ws.Range("E21").Value = "  -1.77%  "

# Row 22
Set-TextValue $ws.Range("D22") "408.90"
$ws.Range("E22").Value = "  -1.97%  "

# Row 23
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D23") "4.02"
$ws.Range("E23").Value = "  -0.38%  "

# Row 24
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D24") "4.37"
$ws.Range("E24").Value = "  +2.06%  "

# Row 25
Set-TextValue $ws.Range("D25") "13.73"
$ws.Range("E25").Value = "  +8.70%  "

# Row 26
Set-TextValue $ws.Range("D26") "82.78"
$ws.Range("E26").Value = "  -4.24%  "

# Row 27
Set-TextValue $ws.Range("D27") "10.49"
$ws.Range("E27").Value = "  -3.13%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.71"
$ws.Range("E28").Value = "  -5.19%  "

# Row 29
Set-TextValue $ws.Range("D29") "8.57"
$ws.Range("E29").Value = "  -6.24%  "

# Row 30
Set-TextValue $ws.Range("D30") "28.91"
$ws.Range("E30").Value = "  -4.52%  "

# Row 31
Set-TextValue $ws.Range("D31") "6.37"
$ws.Range("E31").Value = "  -2.97%  "

# Row 32
Set-TextValue $ws.Range("D32") "578.47"
$ws.Range("E32").Value = "  -5.35%  "

# Row 33
Set-TextValue $ws.Range("D33") "11.31"
$ws.Range("E33").Value = "  -3.81%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.106"
$ws.Range("E34").Value = "  -3.61%  "

# Row 35
Set-TextValue $ws.Range("D35") "57.69"
$ws.Range("E35").Value = "  -2.69%  "

# Row 36
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D36") "1.00"
$ws.Range("E36").Value = "  -0.08%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D37") "0.147"
$ws.Range("E37").Value = "  +1.23%  "

# Row 38
Set-TextValue $ws.Range("D38") "34.96"
$ws.Range("E38").Value = "  -6.90%  "

# Row 39
Set-TextValue $ws.Range("D39") "3.40"
$ws.Range("E39").Value = "  +2.99%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0₃0732"
$ws.Range("E40").Value = "  -7.45%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.364"
$ws.Range("E41").Value = "  -4.55%  "

# Row 42
Set-TextValue $ws.Range("D42") "3.110.24"
$ws.Range("E42").Value = "  -5.04%  "

# Row 43
Set-TextValue $ws.Range("D43") "1.00"
$ws.Range("E43").Value = "  +0.69%  "

# Row 44
Set-TextValue $ws.Range("D44") "2.78"
$ws.Range("E44").Value = "  -1.12%  "

# Row 45
Set-TextValue $ws.Range("D45") "3.21"
$ws.Range("E45").Value = "  -2.30%  "

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D46") "0.0398"
$ws.Range("E46").Value = "  -3.68%  "

# Row 47
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D47") "2.40"
$ws.Range("E47").Value = "  -5.75%  "

# Row 48
$ws.Range("E48").Value = "  -3.85%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.127"
$ws.Range("E49").Value = "  -4.01%  "

# Row 50
Set-TextValue $ws.Range("D50") "132.47"
$ws.Range("E50").Value = "  -3.73%  "

# Row 51
Set-TextValue $ws.Range("D51") "7.98"
$ws.Range("E51").Value = "  -4.95%  "
